# Apply attendance counts: set specific cells from 0 to 1 in Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3: G3, H3
$ws.Range("G3").Value = 1
$ws.Range("H3").Value = 1

# Row 4: D4, E4
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 1

# Row 5: D5, E5
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 1

# Row 6: D6, E6
$ws.Range("D6").Value = 1
$ws.Range("E6").Value = 1

# Row 7: H7
$ws.Range("H7").Value = 1

# Row 8: H8
$ws.Range("H8").Value = 1

# Row 9: D9, E9
$ws.Range("D9").Value = 1
$ws.Range("E9").Value = 1

# Row 10: H10
$ws.Range("H10").Value = 1

# Row 11: H11
$ws.Range("H11").Value = 1

# Row 12: H12
$ws.Range("H12").Value = 1

# Row 13: H13
$ws.Range("H13").Value = 1

# Row 14: H14
$ws.Range("H14").Value = 1

# Row 15: H15
$ws.Range("H15").Value = 1

# Row 16: H16
$ws.Range("H16").Value = 1

# Row 17: H17
$ws.Range("H17").Value = 1

# Row 18: H18
$ws.Range("H18").Value = 1
